# Apply the "Update CDA Logical model for ST.r2b" edit:
#  - Rename "Include from NullFlavor" sheet to "Include #0"
#  - Bump the Version and Date values on the Metadata sheet
#  - Insert a new "Jurisdiction" property row on the Metadata sheet

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from NullFlavor")

# 1. Rename the include sheet.
$wsInclude.Name = "Include #0"

# 2. Update the Version value (row 3, column B).
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 3. Update the Date value (row 8, column B).
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 4. Insert a new row above row 11 ("Description") for "Jurisdiction",
#    copying the formatting of the preceding data row so the new cells
#    match the sheet's existing row style.
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
